{"js": "// The document contains a single 20-row x 5-column table of simple math\n// expressions (\"a+b=c\" / \"a-b=c\"). The edit replaces every cell's text\n// with a new expression while leaving all formatting (fonts, sizes,\n// paragraph alignment, table structure) untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, row-major, in the same 20x5 layout as the existing table.\nconst newValues = [\n  [\"40-21=19\", \"77+7=84\", \"64-39=25\", \"3+39=42\", \"65+16=81\"],\n  [\"25-19=6\", \"88+9=97\", \"24+47=71\", \"79+16=95\", \"54-5=49\"],\n  [\"95-56=39\", \"50-19=31\", \"48+6=54\", \"15+39=54\", \"44+18=62\"],\n  [\"81-14=67\", \"21-9=12\", \"6+37=43\", \"86-68=18\", \"60-53=7\"],\n  [\"73-7=66\", \"13+59=72\", \"35+39=74\", \"47-38=9\", \"5+56=61\"],\n  [\"44-37=7\", \"80-56=24\", \"90-8=82\", \"26+48=74\", \"59+18=77\"],\n  [\"43-8=35\", \"29+45=74\", \"36+7=43\", \"10-9=1\", \"29+57=86\"],\n  [\"52-35=17\", \"88+4=92\", \"58+33=91\", \"61-57=4\", \"56+8=64\"],\n  [\"4+7=11\", \"34+37=71\", \"54-16=38\", \"45-17=28\", \"45+19=64\"],\n  [\"69+2=71\", \"19+69=88\", \"29+15=44\", \"60-37=23\", \"82-64=18\"],\n  [\"77+5=82\", \"39+33=72\", \"85+8=93\", \"92-33=59\", \"90-11=79\"],\n  [\"90-34=56\", \"70-25=45\", \"81-68=13\", \"26+38=64\", \"6+18=24\"],\n  [\"52-43=9\", \"68+3=71\", \"59+12=71\", \"93-36=57\", \"95-39=56\"],\n  [\"78-39=39\", \"46-27=19\", \"58+23=81\", \"62-54=8\", \"93-86=7\"],\n  [\"50-4=46\", \"51-4=47\", \"29+35=64\", \"38+19=57\", \"70-46=24\"],\n  [\"50-35=15\", \"76-59=17\", \"71-45=26\", \"41-15=26\", \"28+56=84\"],\n  [\"92-87=5\", \"92-73=19\", \"36+57=93\", \"91-16=75\", \"16+48=64\"],\n  [\"61-56=5\", \"5+36=41\", \"51-2=49\", \"8+13=21\", \"3+88=91\"],\n  [\"62-27=35\", \"19+12=31\", \"6+35=41\", \"2+29=31\", \"71-69=2\"],\n  [\"62-29=33\", \"8+24=32\", \"44+39=83\", \"45-17=28\", \"9+48=57\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple math\n# expressions (\"a+b=c\" / \"a-b=c\"). The edit replaces every cell's text\n# with a new expression while leaving all formatting (fonts, sizes,\n# paragraph alignment, table structure) untouched.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# New values, row-major, in the same 20x5 layout as the existing table.\n$newValues = @(\n    @(\"40-21=19\",\"77+7=84\",\"64-39=25\",\"3+39=42\",\"65+16=81\"),\n    @(\"25-19=6\",\"88+9=97\",\"24+47=71\",\"79+16=95\",\"54-5=49\"),\n    @(\"95-56=39\",\"50-19=31\",\"48+6=54\",\"15+39=54\",\"44+18=62\"),\n    @(\"81-14=67\",\"21-9=12\",\"6+37=43\",\"86-68=18\",\"60-53=7\"),\n    @(\"73-7=66\",\"13+59=72\",\"35+39=74\",\"47-38=9\",\"5+56=61\"),\n    @(\"44-37=7\",\"80-56=24\",\"90-8=82\",\"26+48=74\",\"59+18=77\"),\n    @(\"43-8=35\",\"29+45=74\",\"36+7=43\",\"10-9=1\",\"29+57=86\"),\n    @(\"52-35=17\",\"88+4=92\",\"58+33=91\",\"61-57=4\",\"56+8=64\"),\n    @(\"4+7=11\",\"34+37=71\",\"54-16=38\",\"45-17=28\",\"45+19=64\"),\n    @(\"69+2=71\",\"19+69=88\",\"29+15=44\",\"60-37=23\",\"82-64=18\"),\n    @(\"77+5=82\",\"39+33=72\",\"85+8=93\",\"92-33=59\",\"90-11=79\"),\n    @(\"90-34=56\",\"70-25=45\",\"81-68=13\",\"26+38=64\",\"6+18=24\"),\n    @(\"52-43=9\",\"68+3=71\",\"59+12=71\",\"93-36=57\",\"95-39=56\"),\n    @(\"78-39=39\",\"46-27=19\",\"58+23=81\",\"62-54=8\",\"93-86=7\"),\n    @(\"50-4=46\",\"51-4=47\",\"29+35=64\",\"38+19=57\",\"70-46=24\"),\n    @(\"50-35=15\",\"76-59=17\",\"71-45=26\",\"41-15=26\",\"28+56=84\"),\n    @(\"92-87=5\",\"92-73=19\",\"36+57=93\",\"91-16=75\",\"16+48=64\"),\n    @(\"61-56=5\",\"5+36=41\",\"51-2=49\",\"8+13=21\",\"3+88=91\"),\n    @(\"62-27=35\",\"19+12=31\",\"6+35=41\",\"2+29=31\",\"71-69=2\"),\n    @(\"62-29=33\",\"8+24=32\",\"44+39=83\",\"45-17=28\",\"9+48=57\")\n)\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $newValues[$r-1][$c-1]\n    }\n}\n"}
